# Adding "Area" / "Atotal" columns (G, H) to the discharge sheet.
#
# Source data: x / depth / velocity / - / segment / Q / Qtotal (A..F)
# New data:    Area (G) = per-segment wetted area, Atotal (H) = running total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- New headers in row 1 --------------------------------------------------
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# -- D column (midpoint depth) becomes a shared formula over D3:D8 --------
# (same formula text that was already there for each row; re-entering it as
#  one range assignment lets Excel collapse it into a shared formula, which
#  is what the target workbook does)
$ws.Range("D3:D8").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# -- G column: incremental wetted area per station -------------------------
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G11").Formula = "=(D4-D3)*B4/100"

# -- H2: running total of the Area column -----------------------------------
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# -- Match the saved selection shown in the edited workbook -----------------
$ws.Range("H2").Select() | Out-Null
